$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.618.93'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.847.88'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = "'262.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.28%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = "'0.5330"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.83%  '
$ws.Range("D8").Value = "'0.3165"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.99%  '
$ws.Range("D9").Value = "'0.06971"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("E10").Value = '  -0.16%  '
$ws.Range("D11").Value = "'0.7727"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("D12").Value = "'0.07842"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").Value = '1.876.61'
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("D14").Value = "'89.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = "'14.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = "'0.000007988"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '26.634.63'
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.085.72'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = "'4.647"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = "'6.028"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = "'9.372"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.18%  '
$ws.Range("D25").Value = "'2.214"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = "'142.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.76%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = "'1.702"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.00%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'17.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.50%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = "'111.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'4.318"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.66%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = "'0.08779"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'4.114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.04866"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = "'0.7405"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.45%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'1.140"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = "'2.888"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.28%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = "'3.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = "'2.358"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.47%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.01741"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.63%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.4827"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'0.9063"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = "'108.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = "'5.915"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.88%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = "'1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = "'7.718"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'0.4207"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'9.127"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = "'0.1251"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = "'35.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.05816"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.32%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = "'0.8986"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.82%  '
